$wb = $excel.ActiveWorkbook

# Rename the metadata sheet
$ws = $wb.Worksheets.Item("SwateTemplateMetadata")
$ws.Name = "isa_template"

# Row 12 "Tags": drop the ER-duplicate "PRIDE" tag, shifting remaining tags left
$ws.Range("B12").Value = "Proteomics"
$ws.Range("C12").Value = "Sample"
$ws.Range("D12").Value = "Mass spectrometry"
$ws.Range("E12").Value = "MS"
$ws.Range("F12").Value = "Preparation"
$ws.Range("G12").Clear()

# Row 13 "Tags Term Accession Number": drop the ER-duplicate accession, shift left
$ws.Range("B13").Value = "http://purl.obolibrary.org/obo/NCIT_C20085"
$ws.Range("C13").ClearContents()
$ws.Range("G13").Clear()

# Row 14 "Tags Term Source REF": drop the ER-duplicate source ref, shift left
$ws.Range("B14").Value = "NCIT"
$ws.Range("C14").ClearContents()
$ws.Range("G14").Clear()

# Row heights follow the content that now wraps in row 12 / no longer wraps in row 13
$ws.Rows.Item(12).RowHeight = 57.6
$ws.Rows.Item(13).RowHeight = 14.4

$ws.Range("B12").Select()
